$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prices (column D) are stored as plain text in this workbook (e.g. thousand
# separators using '.', trailing zeros that must be preserved, etc.), so we
# prefix numeric-looking values with an apostrophe to force text entry and
# avoid Excel's automatic number coercion (which would strip trailing zeros
# or switch to scientific notation).

$ws.Range("D2").Value = "'26.913.59"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "'1.816.05"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'309.14"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "'0.4689"
$ws.Range("E7").Value = "  +1.23%  "

$ws.Range("E8").Value = "  -1.67%  "

$ws.Range("D9").Value = "'0.07380"
$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("D10").Value = "'0.8712"
$ws.Range("E10").Value = "  -0.84%  "

$ws.Range("E11").Value = "  -0.36%  "

$ws.Range("D12").Value = "'1.862.33"
$ws.Range("E12").Value = "  +2.57%  "

$ws.Range("D13").Value = "'5.387"
$ws.Range("E13").Value = "  +0.50%  "

$ws.Range("D14").Value = "'6.525"
$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").Value = "'0.07071"
$ws.Range("E15").Value = "  +0.30%  "

$ws.Range("D16").Value = "'91.73"
$ws.Range("E16").Value = "  +0.26%  "

$ws.Range("E17").Value = "  +0.20%  "

$ws.Range("D18").Value = "'0.000008722"
$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").Value = "'26.959.81"
$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("D22").Value = "'5.323"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("E23").Value = "  -2.33%  "

$ws.Range("D24").Value = "'2.065.09"
$ws.Range("E24").Value = "  +1.91%  "

$ws.Range("D25").Value = "'1.891"
$ws.Range("E25").Value = "  -1.70%  "

$ws.Range("D26").Value = "'150.92"
$ws.Range("E26").Value = "  -0.40%  "

$ws.Range("D27").Value = "'2.181"
$ws.Range("E27").Value = "  +1.22%  "

$ws.Range("D28").Value = "'18.41"
$ws.Range("E28").Value = "  -1.30%  "

$ws.Range("D29").Value = "'5.346"
$ws.Range("E29").Value = "  +0.72%  "

$ws.Range("D30").Value = "'116.26"
$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("D31").Value = "'0.08962"
$ws.Range("E31").Value = "  +0.61%  "

$ws.Range("D32").Value = "'0.7697"
$ws.Range("E32").Value = "  -0.22%  "

$ws.Range("D33").Value = "'1.165"
$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("D34").Value = "'4.514"
$ws.Range("E34").Value = "  +0.57%  "

$ws.Range("D35").Value = "'2.913"
$ws.Range("E35").Value = "  +0.45%  "

$ws.Range("D36").Value = "'1.002"
$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").Value = "'1.086"
$ws.Range("E37").Value = "  -2.89%  "

# Rows 38 and 39 swap places: VeChain <-> Hedera
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05302"
$ws.Range("E38").Value = "  +1.22%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01965"
$ws.Range("E39").Value = "  +0.30%  "

$ws.Range("D40").Value = "'2.952"
$ws.Range("E40").Value = "  +1.50%  "

$ws.Range("E41").Value = "  +0.15%  "

$ws.Range("D42").Value = "'0.5338"
$ws.Range("E42").Value = "  -0.39%  "

$ws.Range("D43").Value = "'2.354"
$ws.Range("E43").Value = "  -3.46%  "

$ws.Range("D44").Value = "'0.1659"
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("D45").Value = "'8.465"
$ws.Range("E45").Value = "  -1.34%  "

$ws.Range("D46").Value = "'0.4931"
$ws.Range("E46").Value = "  -2.88%  "

$ws.Range("D47").Value = "'10.44"
$ws.Range("E47").Value = "  +1.00%  "

$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  +0.11%  "

$ws.Range("D49").Value = "'1.674"
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("D50").Value = "'103.78"
$ws.Range("E50").Value = "  -0.58%  "

$ws.Range("D51").Value = "'0.06303"
$ws.Range("E51").Value = "  -0.42%  "
